$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 622.95654
$ws.Range("I33").Value = 448.7143
$ws.Range("J33").Value = 894
$ws.Range("K33").Value = 448.7143
$ws.Range("L33").Value = 894
$ws.Range("M33").Value = -219.7143
$ws.Range("N33").Value = -1352
$ws.Range("H53").Value = 304.44446
$ws.Range("I53").Value = 238.11111
$ws.Range("J53").Value = 337.6111
$ws.Range("K53").Value = 238.11111
$ws.Range("L53").Value = 337.6111
$ws.Range("M53").Value = 398.88889
$ws.Range("N53").Value = -1611.6111
$ws.Range("H62").Value = 3970.4443
$ws.Range("I62").Value = 2997.25
$ws.Range("J62").Value = 4749
$ws.Range("K62").Value = 2997.25
$ws.Range("L62").Value = 4749
$ws.Range("M62").Value = -2373.25
$ws.Range("N62").Value = -5997
$ws.Range("H65").Value = 3970.4443
$ws.Range("I65").Value = 2997.25
$ws.Range("J65").Value = 4749
$ws.Range("K65").Value = 14986.25
$ws.Range("L65").Value = 23745
$ws.Range("M65").Value = -11866.25
$ws.Range("N65").Value = -29985
$ws.Range("H86").Value = 17426.143
$ws.Range("J86").Value = 3760
$ws.Range("L86").Value = 3760
$ws.Range("N86").Value = -6006
$ws.Range("H89").Value = 17426.143
$ws.Range("J89").Value = 3760
$ws.Range("L89").Value = 18800
$ws.Range("N89").Value = -30032
$ws.Range("H92").Value = 864.3
$ws.Range("I92").Value = 806.1429000000001
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 806.1429000000001
$ws.Range("L92").Value = 1000
$ws.Range("M92").Value = 441.8570999999999
$ws.Range("N92").Value = -3496
$ws.Range("H132").Value = 5719804.5
$ws.Range("I132").Value = 6902005
$ws.Range("J132").Value = 5835.3335
$ws.Range("K132").Value = 20706015
$ws.Range("L132").Value = 17506.0005
$ws.Range("M132").Value = -20703485
$ws.Range("N132").Value = -22566.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3447.1428
$ws.Range("I63").Value = 2307.2727
$ws.Range("K63").Value = 2307.2727
$ws.Range("M63").Value = -1621.2727
$ws.Range("H66").Value = 3447.1428
$ws.Range("I66").Value = 2307.2727
$ws.Range("K66").Value = 11536.3635
$ws.Range("M66").Value = -8104.363499999999
$ws.Range("H95").Value = 38104
$ws.Range("J95").Value = 38104
$ws.Range("L95").Value = 38104
$ws.Range("N95").Value = -43596
$ws.Range("H122").Value = 2618.963
$ws.Range("I122").Value = 1782.35
$ws.Range("J122").Value = 5009.2856
$ws.Range("K122").Value = 5347.049999999999
$ws.Range("L122").Value = 15027.8568
$ws.Range("M122").Value = -2897.049999999999
$ws.Range("N122").Value = -19927.8568
$ws.Range("H125").Value = 37939.285
$ws.Range("J125").Value = 37939.285
$ws.Range("L125").Value = 37939.285
$ws.Range("N125").Value = -47779.285
$ws.Range("H132").Value = 23813720
$ws.Range("I132").Value = 32261900
$ws.Range("J132").Value = 5211.4546
$ws.Range("K132").Value = 96785700
$ws.Range("L132").Value = 15634.3638
$ws.Range("M132").Value = -96783170
$ws.Range("N132").Value = -20694.3638
$ws.Range("H139").Value = 43478.332
$ws.Range("J139").Value = 43478.332
$ws.Range("L139").Value = 43478.332
$ws.Range("N139").Value = -53758.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 954.5
$ws.Range("J80").Value = 726.3
$ws.Range("L80").Value = 726.3
$ws.Range("N80").Value = -2722.3
$ws.Range("H83").Value = 954.5
$ws.Range("J83").Value = 726.3
$ws.Range("L83").Value = 3631.5
$ws.Range("N83").Value = -13615.5
$ws.Range("H94").Value = 701.1429000000001
$ws.Range("I94").Value = 701.6
$ws.Range("K94").Value = 701.6
$ws.Range("M94").Value = -250.6
$ws.Range("H134").Value = 4285.864
$ws.Range("I134").Value = 4041.1177
$ws.Range("J134").Value = 5118
$ws.Range("K134").Value = 12123.3531
$ws.Range("L134").Value = 15354
$ws.Range("M134").Value = -9588.3531
$ws.Range("N134").Value = -20424

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2634422.8
$ws.Range("I31").Value = 3573759.5
$ws.Range("J31").Value = 4280
$ws.Range("K31").Value = 3573759.5
$ws.Range("L31").Value = 4280
$ws.Range("M31").Value = -3573464.5
$ws.Range("N31").Value = -4870
$ws.Range("H34").Value = 2634422.8
$ws.Range("I34").Value = 3573759.5
$ws.Range("J34").Value = 4280
$ws.Range("K34").Value = 3573759.5
$ws.Range("L34").Value = 4280
$ws.Range("M34").Value = -3573557.5
$ws.Range("N34").Value = -4684
$ws.Range("H58").Value = 29415308
$ws.Range("I58").Value = 1368.5
$ws.Range("J58").Value = 45459270
$ws.Range("K58").Value = 1368.5
$ws.Range("L58").Value = 45459270
$ws.Range("M58").Value = -1165.5
$ws.Range("N58").Value = -45459676
$ws.Range("H132").Value = 3771.12
$ws.Range("I132").Value = 2737.6667
$ws.Range("J132").Value = 6428.5713
$ws.Range("K132").Value = 8213.000100000001
$ws.Range("L132").Value = 19285.7139
$ws.Range("M132").Value = -5683.000100000001
$ws.Range("N132").Value = -24345.7139
$ws.Range("H134").Value = 2019.625
$ws.Range("I134").Value = 976.5833
$ws.Range("J134").Value = 5148.75
$ws.Range("K134").Value = 2929.7499
$ws.Range("L134").Value = 15446.25
$ws.Range("M134").Value = -394.7498999999998
$ws.Range("N134").Value = -20516.25
$ws.Range("H136").Value = 29415308
$ws.Range("I136").Value = 1368.5
$ws.Range("J136").Value = 45459270
$ws.Range("K136").Value = 4105.5
$ws.Range("L136").Value = 136377810
$ws.Range("M136").Value = -1555.5
$ws.Range("N136").Value = -136382910

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3843.1428
$ws.Range("J80").Value = 4317
$ws.Range("L80").Value = 12951
$ws.Range("N80").Value = -14823
$ws.Range("H83").Value = 3843.1428
$ws.Range("J83").Value = 4317
$ws.Range("L83").Value = 38853
$ws.Range("N83").Value = -48213
$ws.Range("H86").Value = 1045.3334
$ws.Range("I86").Value = 948
$ws.Range("J86").Value = 1240
$ws.Range("K86").Value = 2844
$ws.Range("L86").Value = 3720
$ws.Range("M86").Value = -1658
$ws.Range("N86").Value = -6092
$ws.Range("H89").Value = 1045.3334
$ws.Range("I89").Value = 948
$ws.Range("J89").Value = 1240
$ws.Range("K89").Value = 8532
$ws.Range("L89").Value = 11160
$ws.Range("M89").Value = -2604
$ws.Range("N89").Value = -23016
$ws.Range("H98").Value = 218
$ws.Range("I98").Value = 233.33333
$ws.Range("J98").Value = 195
$ws.Range("K98").Value = 699.99999
$ws.Range("L98").Value = 585
$ws.Range("M98").Value = 798.00001
$ws.Range("N98").Value = -3581
$ws.Range("H131").Value = 64854.703
$ws.Range("J131").Value = 67633.5
$ws.Range("L131").Value = 202900.5
$ws.Range("N131").Value = -212980.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 6561.875
$ws.Range("I5").Value = 3000
$ws.Range("J5").Value = 11141.429
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 11141.429
$ws.Range("M5").Value = -2888
$ws.Range("N5").Value = -11365.429
$ws.Range("H24").Value = 20751.75
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 20751.75
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 20751.75
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -21097.75
$ws.Range("H126").Value = 2978.75
$ws.Range("I126").Value = 1644
$ws.Range("J126").Value = 3312.4375
$ws.Range("K126").Value = 4932
$ws.Range("L126").Value = 9937.3125
$ws.Range("M126").Value = -2462
$ws.Range("N126").Value = -14877.3125
$ws.Range("H132").Value = 2765.1628
$ws.Range("I132").Value = 2100.1724
$ws.Range("K132").Value = 6300.5172
$ws.Range("M132").Value = -3770.5172

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3300
$ws.Range("I7").Value = 1400
$ws.Range("J7").Value = 3775
$ws.Range("K7").Value = 1400
$ws.Range("L7").Value = 3775
$ws.Range("M7").Value = -1288
$ws.Range("N7").Value = -3999
$ws.Range("H14").Value = 432504
$ws.Range("J14").Value = 21255
$ws.Range("L14").Value = 21255
$ws.Range("N14").Value = -21599
$ws.Range("H24").Value = 36503.5
$ws.Range("I24").Value = 3000
$ws.Range("K24").Value = 3000
$ws.Range("M24").Value = -2657
$ws.Range("H122").Value = 2850.7437
$ws.Range("I122").Value = 2442.6333
$ws.Range("K122").Value = 7327.8999
$ws.Range("M122").Value = -4877.8999
$ws.Range("H126").Value = 3300
$ws.Range("I126").Value = 1400
$ws.Range("J126").Value = 3775
$ws.Range("K126").Value = 4200
$ws.Range("L126").Value = 11325
$ws.Range("M126").Value = -1730
$ws.Range("N126").Value = -16265
$ws.Range("H127").Value = 31693.666
$ws.Range("J127").Value = 31693.666
$ws.Range("L127").Value = 31693.666
$ws.Range("N127").Value = -41613.666
$ws.Range("H132").Value = 3933.4285
$ws.Range("I132").Value = 2333.1667
$ws.Range("J132").Value = 4573.533
$ws.Range("K132").Value = 6999.500100000001
$ws.Range("L132").Value = 13720.599
$ws.Range("M132").Value = -4469.500100000001
$ws.Range("N132").Value = -18780.599
$ws.Range("H136").Value = 3849101
$ws.Range("I136").Value = 5558102
$ws.Range("J136").Value = 3848.75
$ws.Range("K136").Value = 16674306
$ws.Range("L136").Value = 11546.25
$ws.Range("M136").Value = -16671756
$ws.Range("N136").Value = -16646.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 18027
$ws.Range("J54").Value = 18027
$ws.Range("L54").Value = 18027
$ws.Range("N54").Value = -19067
$ws.Range("H122").Value = 304699.8
$ws.Range("I122").Value = 346259
$ws.Range("J122").Value = 3395.75
$ws.Range("K122").Value = 1038777
$ws.Range("L122").Value = 10187.25
$ws.Range("M122").Value = -1036327
$ws.Range("N122").Value = -15087.25
$ws.Range("H123").Value = 21450.666
$ws.Range("J123").Value = 21450.666
$ws.Range("L123").Value = 21450.666
$ws.Range("N123").Value = -31250.666
